$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Timestamp column (A2:A97): shift date serials forward by 1 day
$ws.Range("A2").Value = 45919.01041666666
$ws.Range("A3").Value = 45919.02083333334
$ws.Range("A4").Value = 45919.03125
$ws.Range("A5").Value = 45919.04166666666
$ws.Range("A6").Value = 45919.05208333334
$ws.Range("A7").Value = 45919.0625
$ws.Range("A8").Value = 45919.07291666666
$ws.Range("A9").Value = 45919.08333333334
$ws.Range("A10").Value = 45919.09375
$ws.Range("A11").Value = 45919.10416666666
$ws.Range("A12").Value = 45919.11458333334
$ws.Range("A13").Value = 45919.125
$ws.Range("A14").Value = 45919.13541666666
$ws.Range("A15").Value = 45919.14583333334
$ws.Range("A16").Value = 45919.15625
$ws.Range("A17").Value = 45919.16666666666
$ws.Range("A18").Value = 45919.17708333334
$ws.Range("A19").Value = 45919.1875
$ws.Range("A20").Value = 45919.19791666666
$ws.Range("A21").Value = 45919.20833333334
$ws.Range("A22").Value = 45919.21875
$ws.Range("A23").Value = 45919.22916666666
$ws.Range("A24").Value = 45919.23958333334
$ws.Range("A25").Value = 45919.25
$ws.Range("A26").Value = 45919.26041666666
$ws.Range("A27").Value = 45919.27083333334
$ws.Range("A28").Value = 45919.28125
$ws.Range("A29").Value = 45919.29166666666
$ws.Range("A30").Value = 45919.30208333334
$ws.Range("A31").Value = 45919.3125
$ws.Range("A32").Value = 45919.32291666666
$ws.Range("A33").Value = 45919.33333333334
$ws.Range("A34").Value = 45919.34375
$ws.Range("A35").Value = 45919.35416666666
$ws.Range("A36").Value = 45919.36458333334
$ws.Range("A37").Value = 45919.375
$ws.Range("A38").Value = 45919.38541666666
$ws.Range("A39").Value = 45919.39583333334
$ws.Range("A40").Value = 45919.40625
$ws.Range("A41").Value = 45919.41666666666
$ws.Range("A42").Value = 45919.42708333334
$ws.Range("A43").Value = 45919.4375
$ws.Range("A44").Value = 45919.44791666666
$ws.Range("A45").Value = 45919.45833333334
$ws.Range("A46").Value = 45919.46875
$ws.Range("A47").Value = 45919.47916666666
$ws.Range("A48").Value = 45919.48958333334
$ws.Range("A49").Value = 45919.5
$ws.Range("A50").Value = 45919.51041666666
$ws.Range("A51").Value = 45919.52083333334
$ws.Range("A52").Value = 45919.53125
$ws.Range("A53").Value = 45919.54166666666
$ws.Range("A54").Value = 45919.55208333334
$ws.Range("A55").Value = 45919.5625
$ws.Range("A56").Value = 45919.57291666666
$ws.Range("A57").Value = 45919.58333333334
$ws.Range("A58").Value = 45919.59375
$ws.Range("A59").Value = 45919.60416666666
$ws.Range("A60").Value = 45919.61458333334
$ws.Range("A61").Value = 45919.625
$ws.Range("A62").Value = 45919.63541666666
$ws.Range("A63").Value = 45919.64583333334
$ws.Range("A64").Value = 45919.65625
$ws.Range("A65").Value = 45919.66666666666
$ws.Range("A66").Value = 45919.67708333334
$ws.Range("A67").Value = 45919.6875
$ws.Range("A68").Value = 45919.69791666666
$ws.Range("A69").Value = 45919.70833333334
$ws.Range("A70").Value = 45919.71875
$ws.Range("A71").Value = 45919.72916666666
$ws.Range("A72").Value = 45919.73958333334
$ws.Range("A73").Value = 45919.75
$ws.Range("A74").Value = 45919.76041666666
$ws.Range("A75").Value = 45919.77083333334
$ws.Range("A76").Value = 45919.78125
$ws.Range("A77").Value = 45919.79166666666
$ws.Range("A78").Value = 45919.80208333334
$ws.Range("A79").Value = 45919.8125
$ws.Range("A80").Value = 45919.82291666666
$ws.Range("A81").Value = 45919.83333333334
$ws.Range("A82").Value = 45919.84375
$ws.Range("A83").Value = 45919.85416666666
$ws.Range("A84").Value = 45919.86458333334
$ws.Range("A85").Value = 45919.875
$ws.Range("A86").Value = 45919.88541666666
$ws.Range("A87").Value = 45919.89583333334
$ws.Range("A88").Value = 45919.90625
$ws.Range("A89").Value = 45919.91666666666
$ws.Range("A90").Value = 45919.92708333334
$ws.Range("A91").Value = 45919.9375
$ws.Range("A92").Value = 45919.94791666666
$ws.Range("A93").Value = 45919.95833333334
$ws.Range("A94").Value = 45919.96875
$ws.Range("A95").Value = 45919.97916666666
$ws.Range("A96").Value = 45919.98958333334
$ws.Range("A97").Value = 45920

# Update Notified Production (MW) column (B) retrained-model values
$ws.Range("B22").Value = 14
$ws.Range("B23").Value = 14
$ws.Range("B24").Value = 14
$ws.Range("B25").Value = 14
$ws.Range("B26").Value = 95
$ws.Range("B27").Value = 106
$ws.Range("B28").Value = 124
$ws.Range("B29").Value = 145
$ws.Range("B30").Value = 536
$ws.Range("B31").Value = 578
$ws.Range("B32").Value = 643
$ws.Range("B33").Value = 724
$ws.Range("B34").Value = 1235
$ws.Range("B35").Value = 1303
$ws.Range("B36").Value = 1352
$ws.Range("B37").Value = 1401
$ws.Range("B38").Value = 1734
$ws.Range("B39").Value = 1770
$ws.Range("B40").Value = 1808
$ws.Range("B41").Value = 1844
$ws.Range("B42").Value = 2017
$ws.Range("B43").Value = 2039
$ws.Range("B44").Value = 2059
$ws.Range("B45").Value = 2075
$ws.Range("B46").Value = 2157
$ws.Range("B47").Value = 2164
$ws.Range("B48").Value = 2170
$ws.Range("B49").Value = 2177
$ws.Range("B50").Value = 2144
$ws.Range("B51").Value = 2148
$ws.Range("B52").Value = 2144
$ws.Range("B53").Value = 2136
$ws.Range("B54").Value = 1993
$ws.Range("B55").Value = 1982
$ws.Range("B56").Value = 1965
$ws.Range("B57").Value = 1946
$ws.Range("B58").Value = 1786
$ws.Range("B59").Value = 1759
$ws.Range("B60").Value = 1729
$ws.Range("B61").Value = 1699
$ws.Range("B62").Value = 1455
$ws.Range("B63").Value = 1415
$ws.Range("B64").Value = 1373
$ws.Range("B65").Value = 1332
$ws.Range("B66").Value = 936
$ws.Range("B67").Value = 886
$ws.Range("B68").Value = 755
$ws.Range("B69").Value = 715
$ws.Range("B70").Value = 354
$ws.Range("B71").Value = 296
$ws.Range("B72").Value = 221
$ws.Range("B73").Value = 200
$ws.Range("B74").Value = 67
$ws.Range("B76").Value = 21
$ws.Range("B77").Value = 20
$ws.Range("B78").Value = 11
$ws.Range("B79").Value = 11
$ws.Range("B80").Value = 11
$ws.Range("B81").Value = 11
